$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1448.8431
$ws.Range("I15").Value = 1448.8431
$ws.Range("K15").Value = 4346.5293
$ws.Range("M15").Value = -4177.5293
$ws.Range("H19").Value = 951.0714
$ws.Range("J19").Value = 1014.44446
$ws.Range("L19").Value = 1014.44446
$ws.Range("N19").Value = -1364.44446
$ws.Range("H32").Value = 950.3
$ws.Range("J32").Value = 867
$ws.Range("L32").Value = 867
$ws.Range("N32").Value = -1519
$ws.Range("H33").Value = 1906068.8
$ws.Range("I33").Value = 2354443.8
$ws.Range("K33").Value = 2354443.8
$ws.Range("M33").Value = -2354214.8
$ws.Range("H38").Value = 1881
$ws.Range("I38").Value = 1866.125
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 5598.375
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -5226.375
$ws.Range("N38").Value = -6744
$ws.Range("H44").Value = 1111111
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H48").Value = 3590
$ws.Range("J48").Value = 5416.6665
$ws.Range("L48").Value = 16249.9995
$ws.Range("N48").Value = -16833.9995
$ws.Range("H56").Value = 3590
$ws.Range("J56").Value = 5416.6665
$ws.Range("L56").Value = 16249.9995
$ws.Range("N56").Value = -17317.9995
$ws.Range("H58").Value = 4732.364
$ws.Range("J58").Value = 7034.7144
$ws.Range("L58").Value = 21104.1432
$ws.Range("N58").Value = -21404.1432
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H80").Value = 10310.546
$ws.Range("I80").Value = 15611.286
$ws.Range("J80").Value = 1034.25
$ws.Range("K80").Value = 46833.858
$ws.Range("L80").Value = 3102.75
$ws.Range("M80").Value = -45835.858
$ws.Range("N80").Value = -5098.75
$ws.Range("H83").Value = 10310.546
$ws.Range("I83").Value = 15611.286
$ws.Range("J83").Value = 1034.25
$ws.Range("K83").Value = 140501.574
$ws.Range("L83").Value = 9308.25
$ws.Range("M83").Value = -135509.574
$ws.Range("N83").Value = -19292.25
$ws.Range("H87").Value = 73400
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496
$ws.Range("H90").Value = 73400
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480
$ws.Range("H98").Value = 4499.2
$ws.Range("I98").Value = 3658.4
$ws.Range("J98").Value = 5340
$ws.Range("K98").Value = 3658.4
$ws.Range("L98").Value = 5340
$ws.Range("M98").Value = -2160.4
$ws.Range("N98").Value = -8336
$ws.Range("H107").Value = 63092.875
$ws.Range("I107").Value = 67169.07000000001
$ws.Range("K107").Value = 67169.07000000001
$ws.Range("M107").Value = -65249.07000000001
$ws.Range("H113").Value = 4599.6
$ws.Range("I113").Value = 2999.5
$ws.Range("J113").Value = 4999.625
$ws.Range("K113").Value = 2999.5
$ws.Range("L113").Value = 4999.625
$ws.Range("M113").Value = 254.5
$ws.Range("N113").Value = -11507.625
$ws.Range("J116").Value = 7000
$ws.Range("L116").Value = 7000
$ws.Range("N116").Value = -13884
$ws.Range("H122").Value = 4499.2
$ws.Range("I122").Value = 3658.4
$ws.Range("J122").Value = 5340
$ws.Range("K122").Value = 10975.2
$ws.Range("L122").Value = 16020
$ws.Range("M122").Value = -8525.200000000001
$ws.Range("N122").Value = -20920
$ws.Range("H132").Value = 2891.8484
$ws.Range("I132").Value = 2675.1936
$ws.Range("J132").Value = 6250
$ws.Range("K132").Value = 8025.5808
$ws.Range("L132").Value = 18750
$ws.Range("M132").Value = -5495.5808
$ws.Range("N132").Value = -23810
$ws.Range("H135").Value = 770549.75
$ws.Range("I135").Value = 870684.3
$ws.Range("K135").Value = 7836158.7
$ws.Range("M135").Value = -7833623.7
$ws.Range("H137").Value = 544930.2
$ws.Range("I137").Value = 558485.9
$ws.Range("J137").Value = 532088
$ws.Range("K137").Value = 1675457.7
$ws.Range("L137").Value = 1596264
$ws.Range("M137").Value = -1672907.7
$ws.Range("N137").Value = -1601364
$ws.Range("H138").Value = 4685.7925
$ws.Range("I138").Value = 3062
$ws.Range("J138").Value = 5268.6924
$ws.Range("K138").Value = 9186
$ws.Range("L138").Value = 15806.0772
$ws.Range("M138").Value = -4046
$ws.Range("N138").Value = -26086.0772

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2042.8701
$ws.Range("I32").Value = 2031.3424
$ws.Range("K32").Value = 2031.3424
$ws.Range("M32").Value = -1744.3424
$ws.Range("H61").Value = 2868.2327
$ws.Range("I61").Value = 1302.9656
$ws.Range("K61").Value = 1302.9656
$ws.Range("M61").Value = -1090.9656
$ws.Range("H74").Value = 1127.84
$ws.Range("I74").Value = 940.9091
$ws.Range("J74").Value = 2498.6667
$ws.Range("K74").Value = 940.9091
$ws.Range("L74").Value = 2498.6667
$ws.Range("M74").Value = -66.90909999999997
$ws.Range("N74").Value = -4246.6667
$ws.Range("H77").Value = 1127.84
$ws.Range("I77").Value = 940.9091
$ws.Range("J77").Value = 2498.6667
$ws.Range("K77").Value = 4704.5455
$ws.Range("L77").Value = 12493.3335
$ws.Range("M77").Value = -336.5455000000002
$ws.Range("N77").Value = -21229.3335
$ws.Range("H97").Value = 3139.75
$ws.Range("I97").Value = 3375.8262
$ws.Range("K97").Value = 3375.8262
$ws.Range("M97").Value = -2879.8262
$ws.Range("H102").Value = 1679.6
$ws.Range("I102").Value = 1599.75
$ws.Range("K102").Value = 1599.75
$ws.Range("M102").Value = 22.25
$ws.Range("H110").Value = 1668673.6
$ws.Range("J110").Value = 3010
$ws.Range("L110").Value = 3010
$ws.Range("N110").Value = -7100
$ws.Range("H132").Value = 3209.923
$ws.Range("I132").Value = 1764.0968
$ws.Range("J132").Value = 8812.5
$ws.Range("K132").Value = 5292.2904
$ws.Range("L132").Value = 26437.5
$ws.Range("M132").Value = -2762.2904
$ws.Range("N132").Value = -31497.5
$ws.Range("H136").Value = 2868.2327
$ws.Range("I136").Value = 1302.9656
$ws.Range("K136").Value = 3908.8968
$ws.Range("M136").Value = -1358.8968

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 105149.6
$ws.Range("J58").Value = 105149.6
$ws.Range("L58").Value = 105149.6
$ws.Range("N58").Value = -105737.6
$ws.Range("H105").Value = 57584.945
$ws.Range("I105").Value = 64622.5
$ws.Range("K105").Value = 64622.5
$ws.Range("M105").Value = -62875.5
$ws.Range("H134").Value = 16149.026
$ws.Range("I134").Value = 1925.0984
$ws.Range("J134").Value = 73993
$ws.Range("K134").Value = 5775.2952
$ws.Range("L134").Value = 221979
$ws.Range("M134").Value = -3240.2952
$ws.Range("N134").Value = -227049

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2641.8125
$ws.Range("I16").Value = 2590.7856
$ws.Range("K16").Value = 2590.7856
$ws.Range("M16").Value = -2303.7856
$ws.Range("H31").Value = 503731.1
$ws.Range("I31").Value = 590533.5600000001
$ws.Range("J31").Value = 208602.6
$ws.Range("K31").Value = 590533.5600000001
$ws.Range("L31").Value = 208602.6
$ws.Range("M31").Value = -590238.5600000001
$ws.Range("N31").Value = -209192.6
$ws.Range("H34").Value = 503731.1
$ws.Range("I34").Value = 590533.5600000001
$ws.Range("J34").Value = 208602.6
$ws.Range("K34").Value = 590533.5600000001
$ws.Range("L34").Value = 208602.6
$ws.Range("M34").Value = -590331.5600000001
$ws.Range("N34").Value = -209006.6
$ws.Range("H50").Value = 59948
$ws.Range("J50").Value = 59948
$ws.Range("L50").Value = 59948
$ws.Range("N50").Value = -61198
$ws.Range("H51").Value = 55038.6
$ws.Range("J51").Value = 55038.6
$ws.Range("L51").Value = 55038.6
$ws.Range("N51").Value = -56510.6
$ws.Range("H58").Value = 221985.7
$ws.Range("I58").Value = 402958.28
$ws.Range("J58").Value = 6542.143
$ws.Range("K58").Value = 402958.28
$ws.Range("L58").Value = 6542.143
$ws.Range("M58").Value = -402755.28
$ws.Range("N58").Value = -6948.143
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290
$ws.Range("H60").Value = 39723.75
$ws.Range("I60").Value = 19499.5
$ws.Range("J60").Value = 59948
$ws.Range("K60").Value = 19499.5
$ws.Range("L60").Value = 59948
$ws.Range("M60").Value = -18988.5
$ws.Range("N60").Value = -60970
$ws.Range("H61").Value = 55038.6
$ws.Range("J61").Value = 55038.6
$ws.Range("L61").Value = 55038.6
$ws.Range("N61").Value = -55734.6
$ws.Range("H62").Value = 3568.2856
$ws.Range("I62").Value = 2049.75
$ws.Range("K62").Value = 2049.75
$ws.Range("M62").Value = -1425.75
$ws.Range("H65").Value = 3568.2856
$ws.Range("I65").Value = 2049.75
$ws.Range("K65").Value = 10248.75
$ws.Range("M65").Value = -7128.75
$ws.Range("H86").Value = 10571.143
$ws.Range("J86").Value = 7666.6665
$ws.Range("L86").Value = 7666.6665
$ws.Range("N86").Value = -9912.666499999999
$ws.Range("H89").Value = 10571.143
$ws.Range("J89").Value = 7666.6665
$ws.Range("L89").Value = 38333.3325
$ws.Range("N89").Value = -49565.3325
$ws.Range("H99").Value = 6708.875
$ws.Range("I99").Value = 6574.273
$ws.Range("J99").Value = 6822.769
$ws.Range("K99").Value = 6574.273
$ws.Range("L99").Value = 6822.769
$ws.Range("M99").Value = -5076.273
$ws.Range("N99").Value = -9818.769
$ws.Range("H107").Value = 682.46155
$ws.Range("I107").Value = 608.6923
$ws.Range("K107").Value = 608.6923
$ws.Range("M107").Value = 1311.3077
$ws.Range("H112").Value = 70000
$ws.Range("J112").Value = 70000
$ws.Range("L112").Value = 70000
$ws.Range("N112").Value = -72954
$ws.Range("H113").Value = 2641.8125
$ws.Range("I113").Value = 2590.7856
$ws.Range("K113").Value = 2590.7856
$ws.Range("M113").Value = -420.7856000000002
$ws.Range("H126").Value = 6708.875
$ws.Range("I126").Value = 6574.273
$ws.Range("J126").Value = 6822.769
$ws.Range("K126").Value = 19722.819
$ws.Range("L126").Value = 20468.307
$ws.Range("M126").Value = -17252.819
$ws.Range("N126").Value = -25408.307
$ws.Range("H132").Value = 2285.3157
$ws.Range("I132").Value = 2064.361
$ws.Range("J132").Value = 6262.5
$ws.Range("K132").Value = 6193.083
$ws.Range("L132").Value = 18787.5
$ws.Range("M132").Value = -3663.083
$ws.Range("N132").Value = -23847.5
$ws.Range("H136").Value = 221985.7
$ws.Range("I136").Value = 402958.28
$ws.Range("J136").Value = 6542.143
$ws.Range("K136").Value = 1208874.84
$ws.Range("L136").Value = 19626.429
$ws.Range("M136").Value = -1206324.84
$ws.Range("N136").Value = -24726.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4330.769
$ws.Range("J2").Value = 18
$ws.Range("L2").Value = 108
$ws.Range("N2").Value = -334
$ws.Range("H5").Value = 900908.5600000001
$ws.Range("I5").Value = 267263.34
$ws.Range("J5").Value = 1112123.6
$ws.Range("K5").Value = 801790.02
$ws.Range("L5").Value = 3336370.8
$ws.Range("M5").Value = -801678.02
$ws.Range("N5").Value = -3336594.8
$ws.Range("H22").Value = 2375
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 2833.6667
$ws.Range("K22").Value = 2997
$ws.Range("L22").Value = 8501.000100000001
$ws.Range("M22").Value = -2828
$ws.Range("N22").Value = -8839.000100000001
$ws.Range("H27").Value = 2375
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 2833.6667
$ws.Range("K27").Value = 2997
$ws.Range("L27").Value = 8501.000100000001
$ws.Range("M27").Value = -2895
$ws.Range("N27").Value = -8705.000100000001
$ws.Range("H33").Value = 2057726.5
$ws.Range("I33").Value = 2469255
$ws.Range("J33").Value = 83.333336
$ws.Range("K33").Value = 14815530
$ws.Range("L33").Value = 500.000016
$ws.Range("M33").Value = -14815247
$ws.Range("N33").Value = -1066.000016
$ws.Range("H40").Value = 123
$ws.Range("I40").Value = 87
$ws.Range("J40").Value = 223.8
$ws.Range("K40").Value = 348
$ws.Range("L40").Value = 895.2
$ws.Range("M40").Value = -279
$ws.Range("N40").Value = -1033.2
$ws.Range("H51").Value = 2301.75
$ws.Range("I51").Value = 2654
$ws.Range("K51").Value = 7962
$ws.Range("M51").Value = -7502
$ws.Range("H68").Value = 1671022.8
$ws.Range("I68").Value = 2003613.2
$ws.Range("J68").Value = 1433458.1
$ws.Range("K68").Value = 6010839.6
$ws.Range("L68").Value = 4300374.300000001
$ws.Range("M68").Value = -6010028.6
$ws.Range("N68").Value = -4301996.300000001
$ws.Range("H71").Value = 1671022.8
$ws.Range("I71").Value = 2003613.2
$ws.Range("J71").Value = 1433458.1
$ws.Range("K71").Value = 18032518.8
$ws.Range("L71").Value = 12901122.9
$ws.Range("M71").Value = -18028462.8
$ws.Range("N71").Value = -12909234.9
$ws.Range("H82").Value = 5206
$ws.Range("I82").Value = 5333.3335
$ws.Range("J82").Value = 5015
$ws.Range("K82").Value = 16000.0005
$ws.Range("L82").Value = 15045
$ws.Range("M82").Value = -15594.0005
$ws.Range("N82").Value = -15857
$ws.Range("H85").Value = 5206
$ws.Range("I85").Value = 5333.3335
$ws.Range("J85").Value = 5015
$ws.Range("K85").Value = 16000.0005
$ws.Range("L85").Value = 15045
$ws.Range("M85").Value = -14596.0005
$ws.Range("N85").Value = -17853
$ws.Range("H87").Value = 15000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 15000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H92").Value = 588944.9
$ws.Range("J92").Value = 881.1
$ws.Range("L92").Value = 2643.3
$ws.Range("N92").Value = -5139.3
$ws.Range("H97").Value = 1405.75
$ws.Range("J97").Value = 1206.6666
$ws.Range("L97").Value = 3619.9998
$ws.Range("N97").Value = -4611.9998
$ws.Range("H105").Value = 12515
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12515
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 37545
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -42787
$ws.Range("H107").Value = 53320.69
$ws.Range("J107").Value = 73712.32000000001
$ws.Range("L107").Value = 221136.96
$ws.Range("N107").Value = -224976.96
$ws.Range("H135").Value = 900908.5600000001
$ws.Range("I135").Value = 267263.34
$ws.Range("J135").Value = 1112123.6
$ws.Range("K135").Value = 2405370.06
$ws.Range("L135").Value = 10009112.4
$ws.Range("M135").Value = -2402835.06
$ws.Range("N135").Value = -10014182.4
$ws.Range("H140").Value = 4330.222
$ws.Range("I140").Value = 1386.8572
$ws.Range("K140").Value = 4160.571599999999
$ws.Range("M140").Value = 1019.428400000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 33859.668
$ws.Range("I46").Value = 3880
$ws.Range("J46").Value = 48849.5
$ws.Range("K46").Value = 3880
$ws.Range("L46").Value = 48849.5
$ws.Range("M46").Value = -3724
$ws.Range("N46").Value = -49161.5
$ws.Range("H57").Value = 40000
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41640
$ws.Range("H80").Value = 4007301.8
$ws.Range("I80").Value = 5007502.5
$ws.Range("J80").Value = 3340501.2
$ws.Range("K80").Value = 5007502.5
$ws.Range("L80").Value = 3340501.2
$ws.Range("M80").Value = -5006504.5
$ws.Range("N80").Value = -3342497.2
$ws.Range("H83").Value = 4007301.8
$ws.Range("I83").Value = 5007502.5
$ws.Range("J83").Value = 3340501.2
$ws.Range("K83").Value = 25037512.5
$ws.Range("L83").Value = 16702506
$ws.Range("M83").Value = -25032520.5
$ws.Range("N83").Value = -16712490
$ws.Range("H97").Value = 1416.375
$ws.Range("J97").Value = 2558.2
$ws.Range("L97").Value = 2558.2
$ws.Range("N97").Value = -3550.2
$ws.Range("H122").Value = 4331
$ws.Range("I122").Value = 2717
$ws.Range("J122").Value = 5358.091
$ws.Range("K122").Value = 8151
$ws.Range("L122").Value = 16074.273
$ws.Range("M122").Value = -5701
$ws.Range("N122").Value = -20974.273
$ws.Range("H132").Value = 878350.25
$ws.Range("I132").Value = 1179438.9
$ws.Range("J132").Value = 147135.14
$ws.Range("K132").Value = 3538316.7
$ws.Range("L132").Value = 441405.42
$ws.Range("M132").Value = -3535786.7
$ws.Range("N132").Value = -446465.42

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 606.4
$ws.Range("I9").Value = 255
$ws.Range("K9").Value = 255
$ws.Range("M9").Value = -31
$ws.Range("H16").Value = 1226.25
$ws.Range("J16").Value = 1331.3334
$ws.Range("L16").Value = 1331.3334
$ws.Range("N16").Value = -1671.3334
$ws.Range("H22").Value = 895.44446
$ws.Range("I22").Value = 864.5
$ws.Range("J22").Value = 920.2
$ws.Range("K22").Value = 864.5
$ws.Range("L22").Value = 920.2
$ws.Range("M22").Value = -569.5
$ws.Range("N22").Value = -1510.2
$ws.Range("H27").Value = 895.44446
$ws.Range("I27").Value = 864.5
$ws.Range("J27").Value = 920.2
$ws.Range("K27").Value = 864.5
$ws.Range("L27").Value = 920.2
$ws.Range("M27").Value = -757.5
$ws.Range("N27").Value = -1134.2
$ws.Range("H40").Value = 82160.30499999999
$ws.Range("I40").Value = 146297.86
$ws.Range("J40").Value = 7333.1665
$ws.Range("K40").Value = 146297.86
$ws.Range("L40").Value = 7333.1665
$ws.Range("M40").Value = -146161.86
$ws.Range("N40").Value = -7605.1665
$ws.Range("H68").Value = 50749.227
$ws.Range("I68").Value = 3678.7144
$ws.Range("K68").Value = 3678.7144
$ws.Range("M68").Value = -2929.7144
$ws.Range("H71").Value = 50749.227
$ws.Range("I71").Value = 3678.7144
$ws.Range("K71").Value = 18393.572
$ws.Range("M71").Value = -14649.572
$ws.Range("H122").Value = 1673500.6
$ws.Range("I122").Value = 1672834
$ws.Range("J122").Value = 1674167.4
$ws.Range("K122").Value = 5018502
$ws.Range("L122").Value = 5022502.199999999
$ws.Range("M122").Value = -5016052
$ws.Range("N122").Value = -5027402.199999999
$ws.Range("H132").Value = 2813.2068
$ws.Range("I132").Value = 2061.4492
$ws.Range("J132").Value = 5694.9443
$ws.Range("K132").Value = 6184.3476
$ws.Range("L132").Value = 17084.8329
$ws.Range("M132").Value = -3654.3476
$ws.Range("N132").Value = -22144.8329
$ws.Range("H136").Value = 362012.4
$ws.Range("I136").Value = 628824.4
$ws.Range("K136").Value = 1886473.2
$ws.Range("M136").Value = -1883923.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H41").Value = 8958.6
$ws.Range("I41").Value = 8788.5
$ws.Range("K41").Value = 8788.5
$ws.Range("M41").Value = -8398.5
$ws.Range("H86").Value = 62250
$ws.Range("J86").Value = 62250
$ws.Range("L86").Value = 62250
$ws.Range("N86").Value = -64496
$ws.Range("H89").Value = 62250
$ws.Range("J89").Value = 62250
$ws.Range("L89").Value = 311250
$ws.Range("N89").Value = -322482
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H113").Value = 600.0345
$ws.Range("I113").Value = 405.21054
$ws.Range("J113").Value = 970.2
$ws.Range("K113").Value = 1215.63162
$ws.Range("L113").Value = 2910.6
$ws.Range("M113").Value = 954.3683800000001
$ws.Range("N113").Value = -7250.6
$ws.Range("H132").Value = 31853.371
$ws.Range("I132").Value = 2079.577
$ws.Range("J132").Value = 117866.555
$ws.Range("K132").Value = 6238.731000000001
$ws.Range("L132").Value = 353599.665
$ws.Range("M132").Value = -3708.731000000001
$ws.Range("N132").Value = -358659.665
